# Update "F" column (想去人数 / interested-count) figures to match the
# refreshed data snapshot referenced by the commit
# "Update gh-pages to output generated at 456a3b4".
#
# The workbook contains 4 sheets:
#   1 - 展览     (Exhibitions)
#   2 - 演出     (Performances)
#   3 - 本地生活 (Local life)
#   4 - 全部类型 (All types - aggregate of the above)
#
# Several events had their interested-count (column F) bumped by a small
# amount; the same events also appear (duplicated rows) in the aggregate
# "全部类型" sheet, so each value needs to be updated in both places.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetPerformance = $wb.Worksheets.Item("演出")
$sheetLocalLife   = $wb.Worksheets.Item("本地生活")
$sheetAllTypes    = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$sheetExhibition.Range("F8").Value  = 1861
$sheetExhibition.Range("F12").Value = 1564
$sheetExhibition.Range("F13").Value = 1564
$sheetExhibition.Range("F22").Value = 6829
$sheetExhibition.Range("F23").Value = 7312
$sheetExhibition.Range("F36").Value = 643
$sheetExhibition.Range("F38").Value = 1344
$sheetExhibition.Range("F44").Value = 113

# --- 演出 (sheet2) ---
$sheetPerformance.Range("F3").Value = 17

# --- 本地生活 (sheet3) ---
$sheetLocalLife.Range("F3").Value = 2491

# --- 全部类型 (sheet4) ---
$sheetAllTypes.Range("F10").Value = 17
$sheetAllTypes.Range("F13").Value = 1861
$sheetAllTypes.Range("F17").Value = 1564
$sheetAllTypes.Range("F18").Value = 1564
$sheetAllTypes.Range("F27").Value = 6829
$sheetAllTypes.Range("F28").Value = 7312
$sheetAllTypes.Range("F38").Value = 643
$sheetAllTypes.Range("F42").Value = 1344
$sheetAllTypes.Range("F47").Value = 113
